# Quarterly indexing esoteric bug-fix operation
# Updates the forecast error table values for rows 2-11 (Q1..Q9 quarters)
# and re-indexes the N (count) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.04020811781310986
$ws.Range("C2").Value = 0.4393778005513609
$ws.Range("D2").Value = 0.4185873790753713
$ws.Range("E2").Value = 0.6469832911871615
$ws.Range("F2").Value = 0.6546397523846404
$ws.Range("G2").Value = 37

# Row 3
$ws.Range("B3").Value = 0.160380481824124
$ws.Range("C3").Value = 0.4446889719503624
$ws.Range("D3").Value = 0.3671874000978024
$ws.Range("E3").Value = 0.6059598997440362
$ws.Range("F3").Value = 0.5926395686446457
$ws.Range("G3").Value = 36

# Row 4
$ws.Range("B4").Value = 0.0501434945971165
$ws.Range("C4").Value = 0.3866527438859454
$ws.Range("D4").Value = 0.2890951176385043
$ws.Range("E4").Value = 0.537675662122161
$ws.Range("F4").Value = 0.5431478556575458
$ws.Range("G4").Value = 35

# Row 5
$ws.Range("B5").Value = 0.1382408338576072
$ws.Range("C5").Value = 0.4158547541469218
$ws.Range("D5").Value = 0.3100941732660757
$ws.Range("E5").Value = 0.5568609999506841
$ws.Range("F5").Value = 0.5475411686222322
$ws.Range("G5").Value = 34

# Row 6
$ws.Range("B6").Value = 0.07560015670999451
$ws.Range("C6").Value = 0.3691745166536125
$ws.Range("D6").Value = 0.2747241333589515
$ws.Range("E6").Value = 0.5241413295657493
$ws.Range("F6").Value = 0.5267022622804913
$ws.Range("G6").Value = 33

# Row 7
$ws.Range("B7").Value = 0.1192491126694416
$ws.Range("C7").Value = 0.4028315912834265
$ws.Range("D7").Value = 0.3066657344777716
$ws.Range("E7").Value = 0.5537740825262335
$ws.Range("F7").Value = 0.5494352606605323
$ws.Range("G7").Value = 32

# Row 8
$ws.Range("B8").Value = 0.06922558777690588
$ws.Range("C8").Value = 0.3586838391898814
$ws.Range("D8").Value = 0.2714442011738066
$ws.Range("E8").Value = 0.5210030721347108
$ws.Range("F8").Value = 0.5249194412572022
$ws.Range("G8").Value = 31

# Row 9
$ws.Range("B9").Value = 0.1213405906861701
$ws.Range("C9").Value = 0.4021618232026132
$ws.Range("D9").Value = 0.3128306089004944
$ws.Range("E9").Value = 0.55931262179616
$ws.Range("F9").Value = 0.555325691904054
$ws.Range("G9").Value = 30

# Row 10
$ws.Range("B10").Value = 0.05688217204266421
$ws.Range("C10").Value = 0.3510187655753705
$ws.Range("D10").Value = 0.2692614985437906
$ws.Range("E10").Value = 0.5189041323248357
$ws.Range("F10").Value = 0.5249065084911203
$ws.Range("G10").Value = 29

# Row 11
$ws.Range("B11").Value = 0.1375569543866749
$ws.Range("C11").Value = 0.3880777724859489
$ws.Range("D11").Value = 0.3131151733558528
$ws.Range("E11").Value = 0.5595669516294299
$ws.Range("F11").Value = 0.5523488971977372
$ws.Range("G11").Value = 28
